$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('AE2').Value = 'maa://25251 (92.5), ***maa://21730 (17.19), ***maa://39501 (20.0), *maa://36675 (60.0)'
$ws.Range('W3').Value = 'maa://27396 (85.07), maa://27484 (95.74), maa://27480 (82.35)'
$ws.Range('O8').Value = 'maa://32931 (86.9), *maa://21916 (61.02), maa://23252 (92.31), **maa://22759 (45.45), maa://37496 (100.0)'
$ws.Range('K9').Value = 'maa://22762 (91.57), *maa://39552 (75.0)'
$ws.Range('O9').Value = 'maa://22736 (80.49)'
$ws.Range('C10').Value = '***maa://25695 (19.08), **maa://32237 (37.84), ***maa://34206 (18.18), ***maa://39951 (17.86), **maa://39243 (33.33)'
$ws.Range('S10').Value = 'maa://27395 (96.1), maa://22755 (87.74), **maa://22756 (40.91), ***maa://21737 (10.61)'
$ws.Range('W11').Value = 'maa://36713 (97.87)'
$ws.Range('W12').Value = 'maa://22753 (91.33), *maa://21485 (76.74), maa://37962 (81.25)'
$ws.Range('AE13').Value = '**maa://22737 (30.6), maa://39883 (89.29), *maa://39885 (73.68)'
$ws.Range('O14').Value = 'maa://23250 (98.49), maa://20107 (87.1), maa://22772 (100.0), **maa://22745 (50.0)'
$ws.Range('R14').Value = "'2"
$ws.Range('S14').Value = 'maa://22521 (94.44), maa://42751 (100.0)'
$ws.Range('C15').Value = '*maa://22743 (76.44), maa://22734 (83.33), *maa://30808 (64.29), ***maa://36048 (12.9)'
$ws.Range('C18').Value = 'maa://24570 (96.69)'
$ws.Range('AA19').Value = '*maa://30709 (61.01), *maa://36668 (52.17)'
$ws.Range('K20').Value = 'maa://41331 (84.09)'
$ws.Range('O20').Value = 'maa://37442 (96.77)'
$ws.Range('AE21').Value = 'maa://22524 (94.35), *maa://22432 (74.55)'
$ws.Range('G22').Value = 'maa://25236 (95.89), **maa://21678 (48.94), **maa://22735 (50.0)'
$ws.Range('K23').Value = 'maa://39756 (92.59), maa://39875 (95.92)'
$ws.Range('W24').Value = 'maa://29988 (86.41), maa://23504 (92.94), **maa://22892 (40.14), *maa://25141 (77.05), maa://36663 (80.7), ***maa://22815 (23.08)'
$ws.Range('G25').Value = '*maa://29063 (75.74), *maa://25311 (74.19), ***maa://22725 (4.84)'
$ws.Range('AA25').Value = 'maa://31215 (84.34), *maa://24516 (79.07), maa://26001 (88.89)'
$ws.Range('AE25').Value = 'maa://20108 (96.15), maa://24621 (96.43), maa://36676 (100.0), maa://22771 (84.62), maa://37772 (100.0)'
$ws.Range('AA26').Value = '*maa://42235 (80.0)'
$ws.Range('C28').Value = 'maa://24465 (90.39), maa://25725 (82.28)'
$ws.Range('W28').Value = 'maa://39929 (87.08), ***maa://39723 (14.71), maa://41749 (85.0)'
$ws.Range('AE29').Value = '*maa://24080 (69.11), ***maa://34960 (8.7)'
$ws.Range('W30').Value = '*maa://39477 (77.78)'
$ws.Range('K31').Value = 'maa://35926 (93.19), maa://36258 (80.52)'
$ws.Range('S32').Value = 'maa://41108 (87.5), maa://41238 (94.74)'
$ws.Range('K35').Value = 'maa://41296 (98.25)'
$ws.Range('AE38').Value = 'maa://36697 (83.7)'
$ws.Range('S44').Value = 'maa://39366 (84.0)'
$ws.Range('G46').Value = 'maa://35931 (92.67)'
